$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with P1 and Q1, matching style of existing header cells ---
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("P1:Q1").VerticalAlignment = -4160     # xlTop
$ws.Range("P1:Q1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("P1:Q1").Borders.Weight = 2            # xlThin

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update data rows 2-25: swap I<->K and M<->O column values, add P and Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column = 2
}
